# Commit: "Fruta / hortaliza, semanal"
# A new weekly observation is inserted above the existing row 265, pushing
# every subsequent record (old rows 265..402) down by one row. The sheet's
# used range grows from A1:R402 to A1:R403.
#
# The freshly inserted row 265 repeats the same market/category metadata as
# the (now shifted) row that follows it, but carries a new date (D) and a
# new volume (J) reading.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 265; everything below (old 265..402)
# shifts down to 266..403, and the sheet dimension grows to A1:R403.
$ws.Rows("265:265").Insert()

# Populate the new row 265 with the new weekly record.
$ws.Range("A265").Value = 3
$ws.Range("B265").Value = "Femacal de La Calera"
$ws.Range("C265").Value = "Coquimbo"
$ws.Range("D265").Value = 45029
$ws.Range("E265").Value = 5
$ws.Range("F265").Value = 100112039
$ws.Range("G265").Value = "Ciboulette"
$ws.Range("H265").Value = "Sin especificar"
$ws.Range("I265").Value = "Primera"
$ws.Range("J265").Value = 100
$ws.Range("K265").Value = 1500
$ws.Range("L265").Value = 1500
$ws.Range("M265").Value = 1500
$ws.Range("N265").Value = "$/docena de atados"
$ws.Range("O265").Value = "Provincia de Quillota"
$ws.Range("P265").Value = 500
$ws.Range("Q265").Value = 3
$ws.Range("R265").Value = "Hortaliza"
